$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append: date-serial (col A), col B, col C, col D
$data = @(
    @(44376,0,9,70.67137809187278),
    @(44377,0,9,70.67137809187278),
    @(44378,0,8,62.81900274833137),
    @(44379,0,6,47.11425206124853),
    @(44380,1,7,54.96662740478995),
    @(44381,0,4,31.40950137416569),
    @(44382,0,1,7.852375343541421),
    @(44383,0,1,7.852375343541421),
    @(44384,0,1,7.852375343541421),
    @(44385,0,1,7.852375343541421),
    @(44386,1,2,15.70475068708284),
    @(44387,0,1,7.852375343541421),
    @(44388,0,1,7.852375343541421),
    @(44389,0,1,7.852375343541421),
    @(44390,1,2,15.70475068708284),
    @(44391,0,2,15.70475068708284),
    @(44392,0,2,15.70475068708284),
    @(44393,0,1,7.852375343541421),
    @(44394,0,1,7.852375343541421),
    @(44395,0,1,7.852375343541421),
    @(44396,0,1,7.852375343541421),
    @(44397,1,1,7.852375343541421),
    @(44398,0,1,7.852375343541421),
    @(44399,0,1,7.852375343541421),
    @(44400,1,2,15.70475068708284),
    @(44401,0,2,15.70475068708284),
    @(44402,1,3,23.55712603062426)
)

# Find the last used row in column A (the existing data ends at row 301)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1
$endRow = $startRow + $data.Count - 1

# Copy the formatting (style, borders, number format) of the last existing
# row down across the full block of new rows in one shot.
$srcRow = $ws.Range("A" + $lastRow + ":D" + $lastRow)
$destBlock = $ws.Range("A" + $startRow + ":D" + $endRow)
$srcRow.Copy() | Out-Null
$destBlock.PasteSpecial(-4122) | Out-Null

# Fill in the values row by row.
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
